$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New composition data (labels shift down one row; a new "Abnormal test
# result" category is inserted at row 2; "Symptom – Eye/Ear" category is
# dropped; all percentages are refreshed per corrected pCO2 extraction).
$data = @(
    @{ Row = 2;  Label = "Abnormal test result";        B = 2.4;                 C = 1.7;  D = 2.2 },
    @{ Row = 3;  Label = "Diseases (patient-stated)";   B = 7.4;                 C = 3.1;  D = 3.5 },
    @{ Row = 4;  Label = "Injuries & adverse effects";  B = 7.8;                 C = 4.4;  D = 4.1 },
    @{ Row = 5;  Label = "Other";                       B = 7.5;                 C = 3.9;  D = 8.9 },
    @{ Row = 6;  Label = "Symptom – Circulatory";       B = 9.9;                 C = 6.4;  D = 9.5 },
    @{ Row = 7;  Label = "Symptom – Digestive";         B = 14.7;                C = 6.8;  D = 6 },
    @{ Row = 8;  Label = "Symptom – General";           B = 6.7;                 C = 4.1;  D = 6 },
    @{ Row = 9;  Label = "Symptom – Genitourinary";     B = 6.1;                 C = 4.2;  D = 3.8 },
    @{ Row = 10; Label = "Symptom – Nervous";           B = 14.5;                C = 12.4; D = 8.199999999999999 },
    @{ Row = 11; Label = "Symptom – Respiratory";       B = 20.5;                C = 51.5; D = 44 },
    @{ Row = 12; Label = "Symptom – Skin/Hair/Nails";   B = 2.5;                 C = 1.5;  D = 3.8 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Label
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}
